$wb = $excel.ActiveWorkbook

# --- Kim (sheet1): shift years +1, drop the trailing extra year row, clear bold ---
$kim = $wb.Worksheets.Item("Kim")
$kim.Rows.Item(39).Delete()
for ($r = 2; $r -le 38; $r++) {
    $kim.Cells.Item($r, 1).Value = 2019 + $r
    $kim.Cells.Item($r, 1).Font.Bold = $false
}

# --- Sam (sheet2): shift years +1, drop the trailing extra year row, clear bold ---
$sam = $wb.Worksheets.Item("Sam")
$sam.Rows.Item(39).Delete()
for ($r = 2; $r -le 38; $r++) {
    $sam.Cells.Item($r, 1).Value = 2019 + $r
    $sam.Cells.Item($r, 1).Font.Bold = $false
}

# --- Update per-sheet selections, and switch the active tab from "Fixed Assets" to "Sam" ---
$kim.Activate()
$kim.Range("B7").Select()

$sam.Activate()
$sam.Range("B5").Select()

Write-Host "done"
